# Recalculated Leve profit/price figures for several crafting-job sheets
# (market price refresh from the scheduled data runner).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17: One for the Road
$ws.Range("H17").Value = 940
$ws.Range("J17").Value = 1161.6666
$ws.Range("L17").Value = 3484.9998
$ws.Range("N17").Value = -3820.9998
# Row 18: You Grow, Girl
$ws.Range("H18").Value = 1403.7693
$ws.Range("I18").Value = 1224.9
$ws.Range("K18").Value = 1224.9
$ws.Range("M18").Value = -940.9000000000001
# Row 31: Hush Little Wailer
$ws.Range("H31").Value = 490001.5
$ws.Range("I31").Value = 490001.5
$ws.Range("K31").Value = 1470004.5
$ws.Range("M31").Value = -1469774.5

$ws = $wb.Worksheets.Item("ARM")
# Row 2: Ain't Got No Ingots
$ws.Range("H2").Value = 796.5333000000001
$ws.Range("I2").Value = 621.8
$ws.Range("J2").Value = 1146
$ws.Range("K2").Value = 621.8
$ws.Range("L2").Value = 1146
$ws.Range("M2").Value = -508.8
$ws.Range("N2").Value = -1372
# Row 116: No Scope
$ws.Range("H116").Value = 796.5333000000001
$ws.Range("I116").Value = 621.8
$ws.Range("J116").Value = 1146
$ws.Range("K116").Value = 621.8
$ws.Range("L116").Value = 1146
$ws.Range("M116").Value = 1672.2
$ws.Range("N116").Value = -5734

$ws = $wb.Worksheets.Item("BSM")
# Row 3: Hells Bells
$ws.Range("H3").Value = 796.5333000000001
$ws.Range("I3").Value = 621.8
$ws.Range("J3").Value = 1146
$ws.Range("K3").Value = 621.8
$ws.Range("L3").Value = 1146
$ws.Range("M3").Value = -507.8
$ws.Range("N3").Value = -1374
# Row 105: Ingot to Wing It
$ws.Range("H105").Value = 201978940
$ws.Range("I105").Value = 201978940
$ws.Range("K105").Value = 201978940
$ws.Range("M105").Value = -201977193
# Row 134: Ruthenium Supremium
$ws.Range("H134").Value = 4235.5454
$ws.Range("I134").Value = 933.5862
$ws.Range("J134").Value = 28174.75
$ws.Range("K134").Value = 2800.7586
$ws.Range("L134").Value = 84524.25
$ws.Range("M134").Value = -265.7586000000001
$ws.Range("N134").Value = -89594.25

$ws = $wb.Worksheets.Item("CRP")
# Row 41: The Lone Bowman
$ws.Range("H41").Value = 24566
$ws.Range("J41").Value = 24566
$ws.Range("L41").Value = 24566
$ws.Range("N41").Value = -25422
# Row 51: Greenstone for Greenhorns
$ws.Range("H51").Value = 23251.8
$ws.Range("J51").Value = 24064.75
$ws.Range("L51").Value = 24064.75
$ws.Range("N51").Value = -25536.75
# Row 59: Bow Down to Magic
$ws.Range("H59").Value = 17013.5
$ws.Range("I59").Value = 3900
$ws.Range("K59").Value = 3900
$ws.Range("M59").Value = -2755
# Row 60: Bowing to Greater Power
$ws.Range("H60").Value = 4000
$ws.Range("I60").Value = 4000
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 4000
$ws.Range("L60").Value = 0
$ws.Range("M60").Value = -3489
$ws.Range("N60").ClearContents()
# Row 61: Incant Now, Think Later
$ws.Range("H61").Value = 23251.8
$ws.Range("J61").Value = 24064.75
$ws.Range("L61").Value = 24064.75
$ws.Range("N61").Value = -24760.75
# Row 74: License to Heal
$ws.Range("H74").Value = 26000
$ws.Range("I74").Value = 19000
$ws.Range("K74").Value = 19000
$ws.Range("M74").Value = -18126
# Row 77: Purified Polyrhythm (L)
$ws.Range("H77").Value = 26000
$ws.Range("I77").Value = 19000
$ws.Range("K77").Value = 57000
$ws.Range("M77").Value = -52632

$ws = $wb.Worksheets.Item("CUL")
# Row 131: The Mountain Steeped
$ws.Range("H131").Value = 13891594
$ws.Range("J131").Value = 2964
$ws.Range("L131").Value = 8892
$ws.Range("N131").Value = -18972

$ws = $wb.Worksheets.Item("GSM")
# Row 39: One Man's Trash
$ws.Range("H39").Value = 15000
$ws.Range("J39").Value = 15000
$ws.Range("L39").Value = 15000
$ws.Range("N39").Value = -16064
# Row 80: Needs More Prayerbell
$ws.Range("H80").Value = 2993.0833
$ws.Range("I80").Value = 2318.6428
$ws.Range("K80").Value = 2318.6428
$ws.Range("M80").Value = -1320.6428
# Row 83: With a Noise That Reaches Heaven (L)
$ws.Range("H83").Value = 2993.0833
$ws.Range("I83").Value = 2318.6428
$ws.Range("K83").Value = 11593.214
$ws.Range("M83").Value = -6601.214
# Row 96: Bracelet for Impact
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 74: Overall, We Blend In
$ws.Range("H74").Value = 26400
$ws.Range("J74").Value = 26400
$ws.Range("L74").Value = 26400
$ws.Range("N74").Value = -28396
# Row 77: Eviction Notice (L)
$ws.Range("H77").Value = 26400
$ws.Range("J77").Value = 26400
$ws.Range("L77").Value = 79200
$ws.Range("N77").Value = -89184
# Row 122: Hell on Leather
$ws.Range("H122").Value = 27783676
$ws.Range("I122").Value = 41668216
$ws.Range("K122").Value = 125004648
$ws.Range("M122").Value = -125002198
# Row 132: Tenets of Tanning
$ws.Range("H132").Value = 55452.285
$ws.Range("I132").Value = 18186.428
$ws.Range("J132").Value = 74085.21000000001
$ws.Range("K132").Value = 54559.284
$ws.Range("L132").Value = 222255.63
$ws.Range("M132").Value = -52029.284
$ws.Range("N132").Value = -227315.63

$ws = $wb.Worksheets.Item("WVR")
# Row 8: The Adventurer's New Coat
$ws.Range("H8").Value = 76670.664
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 76670.664
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 76670.664
$ws.Range("M8").ClearContents()
$ws.Range("N8").Value = -76950.664
# Row 11: Wiggle Room
$ws.Range("H11").Value = 3005
$ws.Range("J11").Value = 3005
$ws.Range("L11").Value = 3005
$ws.Range("N11").Value = -3289
# Row 62: Pride Up in Smoke
$ws.Range("H62").Value = 62509376
$ws.Range("I62").Value = 100009000
$ws.Range("J62").Value = 10000.667
$ws.Range("K62").Value = 100009000
$ws.Range("L62").Value = 10000.667
$ws.Range("M62").Value = -100008376
$ws.Range("N62").Value = -11248.667
# Row 65: Desperate for Diversionaries (L)
$ws.Range("H65").Value = 62509376
$ws.Range("I65").Value = 100009000
$ws.Range("J65").Value = 10000.667
$ws.Range("K65").Value = 500045000
$ws.Range("L65").Value = 50003.335
$ws.Range("M65").Value = -500041880
$ws.Range("N65").Value = -56243.335
# Row 81: Where the Dragonflies, the Net Catches
$ws.Range("H81").Value = 634
$ws.Range("I81").Value = 547.5
$ws.Range("J81").Value = 980
$ws.Range("K81").Value = 1095
$ws.Range("L81").Value = 1960
$ws.Range("M81").Value = -34
$ws.Range("N81").Value = -4082
# Row 84: To Kill a Dragon on Nameday (L)
$ws.Range("H84").Value = 634
$ws.Range("I84").Value = 547.5
$ws.Range("J84").Value = 980
$ws.Range("K84").Value = 5475
$ws.Range("L84").Value = 9800
$ws.Range("M84").Value = -171
$ws.Range("N84").Value = -20408
# Row 95: Duress Rehearsal
$ws.Range("H95").Value = 0
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("M95").ClearContents()
$ws.Range("N95").ClearContents()
